$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Mark Du's "Time Spent" (row 6, column B): 34h 17m -> 42h 31m
$ws.Range("B6").Value = "42h 31m"

# Update Stevie Damrel's "Time Spent" (row 8, column B): 24h 30m -> 31h 30m
$ws.Range("B8").Value = "31h 30m"
